# Auto-generated: applies DATA_EXTRACCIO / observation-value refresh
# from the 2026-02-18 19:20 meteocat automatic update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell, far outside the used range, used to stage values that
# look numeric/percent-like ("66%") so a formats-only destination paste
# lands them as literal text without perturbing the destination style.
$scratch = $ws.Range("ZZ1")

$ws.Range("E2").Value = '2026-02-18 19:18:41'
$scratch.Formula = '="66%"'
$scratch.Copy()
$ws.Range("H2").PasteSpecial(-4163)
$ws.Range("E3").Value = '2026-02-18 19:18:44'
$ws.Range("E4").Value = '2026-02-18 19:18:46'
$scratch.Formula = '="75%"'
$scratch.Copy()
$ws.Range("H4").PasteSpecial(-4163)
$ws.Range("J4").Value = '1013.7 hPa'
$ws.Range("L4").Value = '47.9 km/h - 232º 18:44 TU'
$ws.Range("O4").Value = '11.8 °C'
$ws.Range("E5").Value = '2026-02-18 19:18:49'
$ws.Range("E6").Value = '2026-02-18 19:18:52'
$scratch.Formula = '="76%"'
$scratch.Copy()
$ws.Range("H6").PasteSpecial(-4163)
$ws.Range("J6").Value = '1013.3 hPa'
$ws.Range("O6").Value = '11.8 °C'
$ws.Range("E7").Value = '2026-02-18 19:18:55'
$ws.Range("J7").Value = '1014.7 hPa'
$ws.Range("E8").Value = '2026-02-18 19:18:57'
$ws.Range("J8").Value = '1014.5 hPa'
$ws.Range("E9").Value = '2026-02-18 19:19:00'
$scratch.Formula = '="81%"'
$scratch.Copy()
$ws.Range("H9").PasteSpecial(-4163)
$ws.Range("L9").Value = '35.3 km/h - 200º 18:48 TU'
$ws.Range("O9").Value = '10.6 °C'
$ws.Range("E10").Value = '2026-02-18 19:19:03'
$ws.Range("E11").Value = '2026-02-18 19:19:05'
$ws.Range("O11").Value = '5.7 °C'
$ws.Range("E12").Value = '2026-02-18 19:19:08'
$scratch.Formula = '="89%"'
$scratch.Copy()
$ws.Range("H12").PasteSpecial(-4163)
$ws.Range("O12").Value = '10.5 °C'
$ws.Range("E13").Value = '2026-02-18 19:19:10'
$scratch.Formula = '="74%"'
$scratch.Copy()
$ws.Range("H13").PasteSpecial(-4163)
$ws.Range("J13").Value = '1016.0 hPa'
$ws.Range("O13").Value = '3.7 °C'
$ws.Range("E14").Value = '2026-02-18 19:19:13'
$ws.Range("E15").Value = '2026-02-18 19:19:16'
$scratch.Formula = '="80%"'
$scratch.Copy()
$ws.Range("H15").PasteSpecial(-4163)
$ws.Range("O15").Value = '10.5 °C'
$ws.Range("E16").Value = '2026-02-18 19:19:18'
$scratch.Formula = '="49%"'
$scratch.Copy()
$ws.Range("H16").PasteSpecial(-4163)
$ws.Range("E17").Value = '2026-02-18 19:19:21'
$scratch.Formula = '="86%"'
$scratch.Copy()
$ws.Range("H17").PasteSpecial(-4163)
$ws.Range("L17").Value = '108.7 km/h - 259º 18:51 TU'
$ws.Range("E18").Value = '2026-02-18 19:19:24'
$ws.Range("J18").Value = '1013.8 hPa'
$ws.Range("O18").Value = '11.6 °C'
$ws.Range("E19").Value = '2026-02-18 19:19:26'
$scratch.Formula = '="89%"'
$scratch.Copy()
$ws.Range("H19").PasteSpecial(-4163)
$ws.Range("E20").Value = '2026-02-18 19:19:29'
$scratch.Formula = '="73%"'
$scratch.Copy()
$ws.Range("H20").PasteSpecial(-4163)
$ws.Range("E21").Value = '2026-02-18 19:19:32'
$scratch.Formula = '="70%"'
$scratch.Copy()
$ws.Range("H21").PasteSpecial(-4163)
$ws.Range("J21").Value = '1015.3 hPa'
$ws.Range("O21").Value = '6.4 °C'
$ws.Range("E22").Value = '2026-02-18 19:19:35'
$ws.Range("E23").Value = '2026-02-18 19:19:38'
$scratch.Formula = '="53%"'
$scratch.Copy()
$ws.Range("H23").PasteSpecial(-4163)
$ws.Range("O23").Value = '0.3 °C'
$ws.Range("E24").Value = '2026-02-18 19:19:40'
$ws.Range("J24").Value = '1015.3 hPa'
$ws.Range("E25").Value = '2026-02-18 19:19:43'
$scratch.Formula = '="47%"'
$scratch.Copy()
$ws.Range("H25").PasteSpecial(-4163)
$ws.Range("E26").Value = '2026-02-18 19:19:45'
$scratch.Formula = '="70%"'
$scratch.Copy()
$ws.Range("H26").PasteSpecial(-4163)
$ws.Range("J26").Value = '1012.6 hPa'
$ws.Range("E27").Value = '2026-02-18 19:19:48'
$scratch.Formula = '="54%"'
$scratch.Copy()
$ws.Range("H27").PasteSpecial(-4163)
$ws.Range("K27").Value = '13.0 MJ/m2'
$ws.Range("O27").Value = '1.7 °C'
$ws.Range("E28").Value = '2026-02-18 19:19:50'
$scratch.Formula = '="74%"'
$scratch.Copy()
$ws.Range("H28").PasteSpecial(-4163)
$ws.Range("J28").Value = '1013.5 hPa'
$ws.Range("O28").Value = '9.6 °C'
$ws.Range("E29").Value = '2026-02-18 19:19:53'
$ws.Range("E30").Value = '2026-02-18 19:19:56'
$scratch.Formula = '="79%"'
$scratch.Copy()
$ws.Range("H30").PasteSpecial(-4163)
$ws.Range("J30").Value = '1013.1 hPa'
$ws.Range("O30").Value = '10.7 °C'
$ws.Range("E31").Value = '2026-02-18 19:19:58'
$ws.Range("J31").Value = '1011.8 hPa'
$ws.Range("O31").Value = '12.6 °C'
$ws.Range("E32").Value = '2026-02-18 19:20:01'
$scratch.Formula = '="83%"'
$scratch.Copy()
$ws.Range("H32").PasteSpecial(-4163)
$ws.Range("O32").Value = '8.3 °C'
$ws.Range("E33").Value = '2026-02-18 19:20:03'
$ws.Range("J33").Value = '1014.6 hPa'
$ws.Range("K33").Value = '13.6 MJ/m2'
$ws.Range("O33").Value = '4.9 °C'
$ws.Range("E34").Value = '2026-02-18 19:20:06'
$scratch.Formula = '="46%"'
$scratch.Copy()
$ws.Range("H34").PasteSpecial(-4163)
$ws.Range("L34").Value = '51.8 km/h - 48º 18:57 TU'
$ws.Range("E35").Value = '2026-02-18 19:20:09'
$ws.Range("O35").Value = '9.5 °C'
$ws.Range("E36").Value = '2026-02-18 19:20:11'
$ws.Range("J36").Value = '1013.6 hPa'
$ws.Range("L36").Value = '39.2 km/h - 173º 18:48 TU'
$ws.Range("O36").Value = '11.6 °C'
$ws.Range("E37").Value = '2026-02-18 19:20:14'
$ws.Range("J37").Value = '1015.2 hPa'
$ws.Range("O37").Value = '5.9 °C'
$ws.Range("E38").Value = '2026-02-18 19:20:16'
$scratch.Formula = '="76%"'
$scratch.Copy()
$ws.Range("H38").PasteSpecial(-4163)
$ws.Range("O38").Value = '12.4 °C'
$ws.Range("E39").Value = '2026-02-18 19:20:19'
$ws.Range("E40").Value = '2026-02-18 19:20:21'
$scratch.Formula = '="76%"'
$scratch.Copy()
$ws.Range("H40").PasteSpecial(-4163)
$ws.Range("J40").Value = '1016.0 hPa'
$ws.Range("E41").Value = '2026-02-18 19:20:24'
$scratch.Formula = '="85%"'
$scratch.Copy()
$ws.Range("H41").PasteSpecial(-4163)
$ws.Range("J41").Value = '1014.9 hPa'
$ws.Range("E42").Value = '2026-02-18 19:20:26'
$ws.Range("O42").Value = '11.6 °C'
$ws.Range("E43").Value = '2026-02-18 19:20:29'
$scratch.Formula = '="80%"'
$scratch.Copy()
$ws.Range("H43").PasteSpecial(-4163)
$ws.Range("O43").Value = '10.0 °C'
$ws.Range("E44").Value = '2026-02-18 19:20:31'
$scratch.Formula = '="71%"'
$scratch.Copy()
$ws.Range("H44").PasteSpecial(-4163)
$ws.Range("O44").Value = '-1.5 °C'
$ws.Range("E45").Value = '2026-02-18 19:20:34'
$ws.Range("J45").Value = '1012.3 hPa'
$ws.Range("O45").Value = '7.4 °C'
$ws.Range("E46").Value = '2026-02-18 19:20:37'
$scratch.Formula = '="83%"'
$scratch.Copy()
$ws.Range("H46").PasteSpecial(-4163)
$ws.Range("J46").Value = '1015.3 hPa'

$scratch.ClearContents()
$excel.CutCopyMode = $false
